# The source data table restated its 2018-2021 monthly series starting each
# year's row-block with Oct/Nov/Dec (instead of Jan-first), and appended new
# monthly observations through 2022 and 2023-07. Rebuild sheet1 in the
# row order/values the refreshed extract uses.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# date, col B, col C, col D
$data = @(
    @("2018-10", 97.8,               110.6,              98.40000000000001),
    @("2018-11", 95.90000000000001,  113.1,              98.09999999999999),
    @("2018-12", 94.7,               109.3,              99.5),
    @("2018-01", 116.4813,           124.973,            100.3647),
    @("2018-02", 115.7,              128,                100),
    @("2018-03", 113.4,              130.1,              99),
    @("2018-04", 110.7,              127.1,              98.3),
    @("2018-05", 108.8,              125.6,              97.09999999999999),
    @("2018-06", 108.6,              125.4,              96.59999999999999),
    @("2018-07", 106.3,              123.8,              97.09999999999999),
    @("2018-08", 101.4,              119.9,              98.2),
    @("2018-09", 98.8,               110.7,              97.09999999999999),
    @("2019-10", 94.3,               99.3,               117.3),
    @("2019-11", 93.8,               93.09999999999999,  115.8),
    @("2019-12", 93.90000000000001,  91.09999999999999,  114.9),
    @("2019-01", 93.3,               101,                101),
    @("2019-02", 92.8,               98,                 102.2),
    @("2019-03", 93.59999999999999,  95.90000000000001,  102.4),
    @("2019-04", 95.5,               97.40000000000001,  101.7),
    @("2019-05", 95.8,               98.59999999999999,  106.1),
    @("2019-06", 93.5,               97.3,               109.8),
    @("2019-07", 92.5,               96.3,               113.6),
    @("2019-08", 94.3,               96.2,               115.6),
    @("2019-09", 95.2,               96.8,               119),
    @("2020-10", 102.1,              88.5,               113.2),
    @("2020-11", 104.6,              92.7,               112.6),
    @("2020-12", 108.2,              97.2,               111.7),
    @("2020-01", 94.7,               96.2,               114.6),
    @("2020-02", 94,                 98.3,               115.6),
    @("2020-03", 89.59999999999999,  98.8,               115.7),
    @("2020-04", 88.3,               91.09999999999999,  119.4),
    @("2020-05", 90.09999999999999,  87.09999999999999,  119.7),
    @("2020-06", 92.90000000000001,  90.09999999999999,  116.9),
    @("2020-07", 96.09999999999999,  92,                 117.5),
    @("2020-08", 99.7,               89.5,               118.6),
    @("2020-09", 102,                89.7,               114.3),
    @("2021-10", 120.3,              135.7,              91.59999999999999),
    @("2021-11", 119.9,              133.1,              93.8),
    @("2021-12", 116.4,              137.8,              94.8),
    @("2021-01", 111.9,              99.90000000000001,  110.8),
    @("2021-02", 114.2,              101.8,              106.3),
    @("2021-03", 122.2,              105.9,              103),
    @("2021-04", 126.4,              115.1,              99.90000000000001),
    @("2021-05", 128.7,              121.4,              97.40000000000001),
    @("2021-06", 128.8,              124,                96.40000000000001),
    @("2021-07", 125.4,              127.2,              94.2),
    @("2021-08", 121.1,              134.5,              88.90000000000001),
    @("2021-09", 118.9,              137.9,              90.3),
    @("2022-10", 101.5,              115.3,              103.9),
    @("2022-11", 100.5,              118.1,              105),
    @("2022-12", 101.3,              115.9,              106.5),
    @("2022-01", 114.2,              138.3,              94.3),
    @("2022-02", 114.3,              138.3,              97.40000000000001),
    @("2022-03", 112.9,              134,                102.2),
    @("2022-04", 113.2,              133.9,              106.2),
    @("2022-05", 110.9,              132.2,              104.5),
    @("2022-06", 107.8,              123.7,              105.2),
    @("2022-07", 103.8,              120.7,              103.6),
    @("2022-08", 102.5,              115.1,              102.7),
    @("2022-09", 101.7,              112.7,              103.1),
    @("2023-01", 99.8,               119.6,              107.7),
    @("2023-02", 98.7,               126.4,              108.7),
    @("2023-03", 97.8,               128.2,              109.4),
    @("2023-04", 95.2,               114.4,              109.2),
    @("2023-05", 93.5,               103.3,              112.4),
    @("2023-06", 93.3,               113.8,              112.5),
    @("2023-07", 96.8,               117,                115.1)
)

$firstDataRow = 2
$origLastDataRow = 49

$r = $firstDataRow
foreach ($rec in $data) {
    $ws.Cells.Item($r, 1).Value2 = $rec[0]
    $ws.Cells.Item($r, 2).Value2 = $rec[1]
    $ws.Cells.Item($r, 3).Value2 = $rec[2]
    $ws.Cells.Item($r, 4).Value2 = $rec[3]
    $r++
}
$lastDataRow = $r - 1

# The data table grew from 48 to 67 rows. The newly appended rows need the
# same date-column (A) formatting (bold, bordered, centered) that the
# existing rows already carry - copy it down from a formatted cell.
if ($lastDataRow -gt $origLastDataRow) {
    $ws.Cells.Item($firstDataRow, 1).Copy()
    $ws.Range($ws.Cells.Item($origLastDataRow + 1, 1), $ws.Cells.Item($lastDataRow, 1)).PasteSpecial(-4122)
    $excel.CutCopyMode = 0
}
